$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The original row 4 ("회원 / /Photostagram/member/... / GET / 회원
# 로그인 화면 / 홍길동") is being expanded into a full block describing
# all of the member-related API endpoints. Insert 8 extra rows right
# below the existing row 4 so the old row 5 (상품...) ends up at row 13.
# ------------------------------------------------------------------
$ws.Rows("5:12").Insert()

# Clear any formatting the Insert() may have copied down from row 4
# so the new rows match the plain (unstyled) look of the rest of the
# sheet, and drop the now-stale centered style that used to sit on the
# 상품/고객센터/관리자 rows (they are now rows 13-15).
$ws.Range("C5:C12").ClearFormats()
$ws.Range("C4").ClearFormats()
$ws.Range("B13:E15").ClearFormats()

# ---------------------- Row 1 : header (unchanged) -----------------
$ws.Range("A1").Value = "구분"
$ws.Range("B1").Value = "URI"
$ws.Range("C1").Value = "method"
$ws.Range("D1").Value = "설명"
$ws.Range("E1").Value = "작업자"

# ---------------------- Row 2-3 : 메인 (unchanged) ------------------
$ws.Range("A2").Value = "메인"
$ws.Range("B2").Value = "/Photostagram/"
$ws.Range("C2").Value = "GET"
$ws.Range("D2").Value = "메인화면"
$ws.Range("E2").Value = "홍길동"

$ws.Range("B3").Value = "/Photostagram/index.do"
$ws.Range("C3").Value = "GET"
$ws.Range("D3").Value = "메인화면"
$ws.Range("E3").Value = "홍길동"

# ---------------------- Row 4-12 : 회원 (expanded) ------------------
$ws.Range("A4").Value = "회원"
$ws.Range("B4").Value = "/Photostagram/member/login"
$ws.Range("C4").Value = "POST"
$ws.Range("D4").Value = "회원 로그인 화면"
$ws.Range("E4").Value = "김진우"

$ws.Range("B5").Value = "/Photostagram/member/register"
$ws.Range("C5").Value = "POST"
$ws.Range("D5").Value = "회원 가입 화면"
$ws.Range("E5").Value = "김진우"

$ws.Range("B6").Value = "/Photostagram/member/birth"
$ws.Range("C6").Value = "GET"
$ws.Range("D6").Value = "회원가입 생일 입력 화면"
$ws.Range("E6").Value = "김진우"

$ws.Range("B7").Value = "/Photostagram/member/email"
$ws.Range("C7").Value = "POST"
$ws.Range("D7").Value = "회원 가입 이메일 인증 화면"
$ws.Range("E7").Value = "김진우"

$ws.Range("B8").Value = "/Photostagram/member/terms"
$ws.Range("C8").Value = "POST"
$ws.Range("D8").Value = "회원 가입 이용약관 화면"
$ws.Range("E8").Value = "김진우"

$ws.Range("B9").Value = "/Photostagram/member/checkId"
$ws.Range("C9").Value = "POST"
$ws.Range("D9").Value = "아이디 찾기 화면"
$ws.Range("E9").Value = "김진우"

$ws.Range("B10").Value = "/Photostagram/member/checkPass"
$ws.Range("C10").Value = "POST"
$ws.Range("D10").Value = "비밀번호 찾기 화면"
$ws.Range("E10").Value = "김진우"

$ws.Range("B11").Value = "/Photostagram/member/resultId"
$ws.Range("C11").Value = "GET"
$ws.Range("D11").Value = "아이디 찾기 결과 화면"
$ws.Range("E11").Value = "김진우"

$ws.Range("B12").Value = "/Photostagram/member/resultPass"
$ws.Range("C12").Value = "POST"
$ws.Range("D12").Value = "비밀번호 찾기 임시 비밀번호 발송 화면"
$ws.Range("E12").Value = "김진우"

# ---------------------- Row 13 : 상품 (was row 4, shifted) ----------
$ws.Range("A13").Value = "상품"
$ws.Range("B13").Value = "/Photostagram/product/…"
$ws.Range("C13").Value = "GET"
$ws.Range("D13").Value = "상품 목록 화면"
$ws.Range("E13").Value = "홍길동"

# ---------------------- Row 14 : 고객센터 (was row 5, shifted) ------
$ws.Range("A14").Value = "고객센터"
$ws.Range("B14").Value = "/Photostagram/cs/…"
$ws.Range("C14").Value = "GET"
$ws.Range("D14").Value = "고객센터 메인화면"
$ws.Range("E14").Value = "홍길동"

# ---------------------- Row 15 : 관리자 (was row 6, shifted) --------
$ws.Range("A15").Value = "관리자"
$ws.Range("B15").Value = "/Photostagram/admin/…"
$ws.Range("C15").Value = "GET"
$ws.Range("D15").Value = "관리자 메인화면"
$ws.Range("E15").Value = "홍길동"

# Move the active selection the way the authored file shows it.
$ws.Range("C12").Select()
